# Daily attendance processing - 2025-10-28 07:44:06
# Applies the session-analysis refresh: reorders "Recorded By" name lists,
# updates missing/pending counters, updates per-row shortage counts, widens
# the Status column, and flags the newly-overdue session rows (21) for
# groups B2D / B2E / B2F as "Not Recorded" (pink) instead of "Pending" (yellow).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column I ("Status") gets a bit wider to fit "Not Recorded" ---
# (COM ColumnWidth uses "characters" units which differ from the stored
# OOXML width by the fixed ~0.8333 padding for this Calibri 11 workbook,
# so back that out to land on a stored width of exactly 14.)
$ws.Columns.Item(9).ColumnWidth = 13.166666666666666

# --- "Recorded By" (column G) name-order fixes ---
$gUpdates = @{
    "G2"   = "system, backup@backdoor.com, System"
    "G7"   = "admin@admin.com, System"
    "G11"  = "System, dnasr281@gmail.com"
    "G17"  = "System, dnasr281@gmail.com"
    "G29"  = "system, backup@backdoor.com, System"
    "G34"  = "admin@admin.com, System"
    "G38"  = "System, dnasr281@gmail.com"
    "G44"  = "System, dnasr281@gmail.com"
    "G56"  = "system, backup@backdoor.com, System"
    "G61"  = "admin@admin.com, System"
    "G65"  = "System, dnasr281@gmail.com"
    "G71"  = "System, dnasr281@gmail.com"
    "G90"  = "admin@admin.com, dnasr281@gmail.com"
    "G96"  = "System, dnasr281@gmail.com"
    "G97"  = "System, dnasr281@gmail.com"
    "G99"  = "System, dnasr281@gmail.com"
    "G116" = "admin@admin.com, dnasr281@gmail.com"
    "G122" = "System, dnasr281@gmail.com"
    "G123" = "System, dnasr281@gmail.com"
    "G125" = "System, dnasr281@gmail.com"
    "G142" = "admin@admin.com, dnasr281@gmail.com"
    "G148" = "System, dnasr281@gmail.com"
    "G149" = "System, dnasr281@gmail.com"
    "G151" = "System, dnasr281@gmail.com"
}
foreach ($addr in $gUpdates.Keys) {
    $ws.Range($addr).Value = $gUpdates[$addr]
}

# --- Class-statistics counters (Missing / Pending sessions) ---
$ws.Range("L7").Value = 3
$ws.Range("L8").Value = 24

# --- Per-row shortage recompute for rows 18-20 ---
$ws.Range("P18").Value = 1
$ws.Range("Q18").Value = 5
$ws.Range("P19").Value = 1
$ws.Range("Q19").Value = 5
$ws.Range("P20").Value = 1
$ws.Range("Q20").Value = 5

# --- Flag newly-overdue "session 21" rows as Not Recorded (pink) for B2D/B2E/B2F ---
# Build the new fill+font combo (black text on pink) by stacking formats copied
# from existing cells, so the workbook's style table is reused/extended rather
# than duplicated: start from the existing "Pending" style (bold-free black
# text) then swap in the pink fill already used elsewhere in the sheet.
$pinkSource = $ws.Range("L26")      # existing cell formatted with the pink fill
$textSource = $ws.Range("I103")     # existing "Pending" style (black text, centered)

$overdueRows = @(103, 129, 155)
foreach ($r in $overdueRows) {
    $rowRange = $ws.Range("A" + $r + ":I" + $r)

    $textSource.Copy()
    $rowRange.PasteSpecial(-4122)

    $pinkSource.Copy()
    $rowRange.Interior.ColorIndex = $pinkSource.Interior.ColorIndex

    $ws.Range("I" + $r).Value = "Not Recorded"
}

$excel.CutCopyMode = $false
